$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.85
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 2.63
$ws.Range("L3").Value = 5
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("X3").Value = 7.5
$ws.Range("Z3").Value = 15
$ws.Range("AD3").Value = 6.5
$ws.Range("AF3").Value = 81
$ws.Range("AH3").Value = 10
$ws.Range("AJ3").Value = 17
$ws.Range("AK3").Value = 51
$ws.Range("AV3").Value = 81
$ws.Range("AZ3").Value = 101
$ws.Range("BA3").Value = 151

# Row 7
$ws.Range("G7").Value = 1.42
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 9.5
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 2.1
$ws.Range("L7").Value = 8.5
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("S7").Value = 1.5
$ws.Range("T7").Value = 2.5
$ws.Range("U7").Value = 2.63
$ws.Range("V7").Value = 1.44
$ws.Range("W7").Value = 4.75
$ws.Range("X7").Value = 5.5
$ws.Range("Z7").Value = 8.5
$ws.Range("AC7").Value = 7
$ws.Range("AE7").Value = 26
$ws.Range("AF7").Value = 101
$ws.Range("AH7").Value = 17
$ws.Range("AJ7").Value = 29
$ws.Range("AK7").Value = 126
$ws.Range("AL7").Value = 81
$ws.Range("AM7").Value = 81
$ws.Range("AN7").Value = 3.1
$ws.Range("AO7").Value = 7
$ws.Range("AP7").Value = 26
$ws.Range("AS7").Value = 251
$ws.Range("AT7").Value = 2.5
$ws.Range("AV7").Value = 101
$ws.Range("AW7").Value = 9
$ws.Range("AZ7").Value = 251
$ws.Range("BA7").Value = 351

# Row 8
$ws.Range("G8").Value = 1.53
$ws.Range("I8").Value = 6.25
$ws.Range("J8").Value = 2.1
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("Q8").Value = 1.85
$ws.Range("R8").Value = 2
$ws.Range("X8").Value = 7.5
$ws.Range("AG8").Value = 251
$ws.Range("AJ8").Value = 19
$ws.Range("AM8").Value = 41
$ws.Range("AS8").Value = 126
$ws.Range("AX8").Value = 29
$ws.Range("AZ8").Value = 101

# Row 15
$ws.Range("G15").Value = 5.75
$ws.Range("H15").Value = 4.2
$ws.Range("I15").Value = 1.53
$ws.Range("J15").Value = 6
$ws.Range("L15").Value = 2.1
$ws.Range("S15").Value = 1.36
$ws.Range("T15").Value = 3
$ws.Range("X15").Value = 29
$ws.Range("AA15").Value = 41
$ws.Range("AC15").Value = 11
$ws.Range("AE15").Value = 19
$ws.Range("AG15").Value = 351
$ws.Range("AI15").Value = 7
$ws.Range("AN15").Value = 7.5
$ws.Range("AQ15").Value = 126
$ws.Range("AS15").Value = 301
$ws.Range("AT15").Value = 3
$ws.Range("AV15").Value = 51
$ws.Range("AY15").Value = 19
